$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the A column coordinate strings (rows 2-9) and B column values (rows 2-9)
$ws.Range("A2").Value = "53.96692989779158, 25.416455740683055"
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = "56.44296749286285, 24.836864193327976"
$ws.Range("B3").Value = 100

$ws.Range("A4").Value = "54.44296015829089, 22.740469234809616"
$ws.Range("B4").Value = 50

$ws.Range("A5").Value = "53.96692989779158, 23.776334979018685"
$ws.Range("B5").Value = 74

$ws.Range("A6").Value = "56.272170765510346, 21.235997558696436"
$ws.Range("B6").Value = 20

$ws.Range("A7").Value = "55.30851893387957, 26.748283126094723"
$ws.Range("B7").Value = 80

$ws.Range("A8").Value = "56.16245896219404, 25.909725142687375"
$ws.Range("B8").Value = 10

$ws.Range("A9").Value = "54.75601343546629, 25.18349116638527"
$ws.Range("B9").Value = 1

# Add new row 10
$ws.Range("A10").Value = "55.761087088687496, 22.701365275429193"
$ws.Range("B10").Value = 1

# Update selection to D5
$ws.Range("D5").Select()
